# Auto-generated edit script: updates market-price derived columns (H-N)
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR to match refreshed
# price-API data, mirroring a scheduled data-refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 591.4
$ws.Range("I2").Value = 338.75
$ws.Range("K2").Value = 338.75
$ws.Range("M2").Value = -225.75
$ws.Range("H9").Value = 700.9231
$ws.Range("I9").Value = 738.5
$ws.Range("K9").Value = 738.5
$ws.Range("M9").Value = -569.5
$ws.Range("H32").Value = 2432.5
$ws.Range("J32").Value = 2725
$ws.Range("L32").Value = 2725
$ws.Range("N32").Value = -3377
$ws.Range("H33").Value = 296.85715
$ws.Range("I33").Value = 312.3846
$ws.Range("K33").Value = 312.3846
$ws.Range("M33").Value = -83.38459999999998
$ws.Range("H40").Value = 5266.5
$ws.Range("I40").Value = 4332.6665
$ws.Range("J40").Value = 6200.3335
$ws.Range("K40").Value = 4332.6665
$ws.Range("L40").Value = 6200.3335
$ws.Range("M40").Value = -4157.6665
$ws.Range("N40").Value = -6550.3335
$ws.Range("H111").Value = 5799.6665
$ws.Range("I111").Value = 5799.6665
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 17398.9995
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -14331.9995
$ws.Range("N111").ClearContents()
$ws.Range("H138").Value = 1534.4375
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4104.231
$ws.Range("I61").Value = 3108.2258
$ws.Range("J61").Value = 7963.75
$ws.Range("K61").Value = 3108.2258
$ws.Range("L61").Value = 7963.75
$ws.Range("M61").Value = -2896.2258
$ws.Range("N61").Value = -8387.75
$ws.Range("H94").Value = 491999.75
$ws.Range("J94").Value = 491999.75
$ws.Range("L94").Value = 491999.75
$ws.Range("N94").Value = -493801.75
$ws.Range("H132").Value = 3146.9333
$ws.Range("I132").Value = 3146.9333
$ws.Range("K132").Value = 9440.7999
$ws.Range("M132").Value = -6910.7999
$ws.Range("H136").Value = 4104.231
$ws.Range("I136").Value = 3108.2258
$ws.Range("J136").Value = 7963.75
$ws.Range("K136").Value = 9324.6774
$ws.Range("L136").Value = 23891.25
$ws.Range("M136").Value = -6774.6774
$ws.Range("N136").Value = -28991.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6612.909
$ws.Range("I86").Value = 2968
$ws.Range("K86").Value = 2968
$ws.Range("M86").Value = -1845
$ws.Range("H89").Value = 6612.909
$ws.Range("I89").Value = 2968
$ws.Range("K89").Value = 14840
$ws.Range("M89").Value = -9224

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2204.7144
$ws.Range("I31").Value = 1985.25
$ws.Range("J31").Value = 2497.3333
$ws.Range("K31").Value = 1985.25
$ws.Range("L31").Value = 2497.3333
$ws.Range("M31").Value = -1690.25
$ws.Range("N31").Value = -3087.3333
$ws.Range("H34").Value = 2204.7144
$ws.Range("I34").Value = 1985.25
$ws.Range("J34").Value = 2497.3333
$ws.Range("K34").Value = 1985.25
$ws.Range("L34").Value = 2497.3333
$ws.Range("M34").Value = -1783.25
$ws.Range("N34").Value = -2901.3333
$ws.Range("H54").Value = 10320.833
$ws.Range("I54").Value = 6083.25
$ws.Range("K54").Value = 6083.25
$ws.Range("M54").Value = -5425.25
$ws.Range("H58").Value = 3216.5557
$ws.Range("I58").Value = 2368.625
$ws.Range("J58").Value = 10000
$ws.Range("K58").Value = 2368.625
$ws.Range("L58").Value = 10000
$ws.Range("M58").Value = -2165.625
$ws.Range("N58").Value = -10406
$ws.Range("H118").Value = 89999.5
$ws.Range("J118").Value = 89999.5
$ws.Range("L118").Value = 89999.5
$ws.Range("N118").Value = -93313.5
$ws.Range("H129").Value = 94950
$ws.Range("J129").Value = 94950
$ws.Range("L129").Value = 94950
$ws.Range("N129").Value = -104950
$ws.Range("H132").Value = 1222.0667
$ws.Range("I132").Value = 1336.4166
$ws.Range("K132").Value = 4009.2498
$ws.Range("M132").Value = -1479.2498
$ws.Range("H136").Value = 3216.5557
$ws.Range("I136").Value = 2368.625
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 7105.875
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -4555.875
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 876.5
$ws.Range("I46").Value = 876.5
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 2629.5
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -2538.5
$ws.Range("N46").ClearContents()
$ws.Range("H52").Value = 1842.3334
$ws.Range("J52").Value = 1842.3334
$ws.Range("L52").Value = 5527.0002
$ws.Range("N52").Value = -6059.0002
$ws.Range("H56").Value = 11168
$ws.Range("I56").Value = 11168
$ws.Range("K56").Value = 11168
$ws.Range("M56").Value = -10638
$ws.Range("H57").Value = 13999.667
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 13999.667
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 41999.001
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -43117.001
$ws.Range("H131").Value = 616
$ws.Range("I131").Value = 616
$ws.Range("K131").Value = 1848
$ws.Range("M131").Value = 3192
$ws.Range("H132").Value = 1398.7778
$ws.Range("I132").Value = 1447.25
$ws.Range("K132").Value = 13025.25
$ws.Range("M132").Value = -10495.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 980.5
$ws.Range("I13").Value = 461
$ws.Range("J13").Value = 1500
$ws.Range("K13").Value = 461
$ws.Range("L13").Value = 1500
$ws.Range("M13").Value = -322
$ws.Range("N13").Value = -1778
$ws.Range("H17").Value = 15000
$ws.Range("J17").Value = 15000
$ws.Range("L17").Value = 15000
$ws.Range("N17").Value = -15336
$ws.Range("H46").Value = 19885.428
$ws.Range("I46").Value = 4651.25
$ws.Range("J46").Value = 25979.1
$ws.Range("K46").Value = 4651.25
$ws.Range("L46").Value = 25979.1
$ws.Range("M46").Value = -4495.25
$ws.Range("N46").Value = -26291.1
$ws.Range("H80").Value = 2712.1428
$ws.Range("I80").Value = 2595
$ws.Range("J80").Value = 2800
$ws.Range("K80").Value = 2595
$ws.Range("L80").Value = 2800
$ws.Range("M80").Value = -1597
$ws.Range("N80").Value = -4796
$ws.Range("H83").Value = 2712.1428
$ws.Range("I83").Value = 2595
$ws.Range("J83").Value = 2800
$ws.Range("K83").Value = 12975
$ws.Range("L83").Value = 14000
$ws.Range("M83").Value = -7983
$ws.Range("N83").Value = -23984
$ws.Range("H113").Value = 1098.75
$ws.Range("I113").Value = 1031.6666
$ws.Range("K113").Value = 1031.6666
$ws.Range("M113").Value = 1138.3334
$ws.Range("H132").Value = 2543.625
$ws.Range("I132").Value = 2543.625
$ws.Range("K132").Value = 7630.875
$ws.Range("M132").Value = -5100.875
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2972.4614
$ws.Range("J82").Value = 2386.5
$ws.Range("L82").Value = 2386.5
$ws.Range("N82").Value = -3108.5
$ws.Range("H85").Value = 2972.4614
$ws.Range("J85").Value = 2386.5
$ws.Range("L85").Value = 2386.5
$ws.Range("N85").Value = -4882.5
$ws.Range("H118").Value = 44799
$ws.Range("J118").Value = 44799
$ws.Range("L118").Value = 44799
$ws.Range("N118").Value = -48113
$ws.Range("H132").Value = 1988.9286
$ws.Range("I132").Value = 1961.5385
$ws.Range("K132").Value = 5884.6155
$ws.Range("M132").Value = -3354.6155
$ws.Range("H136").Value = 7574.4443
$ws.Range("I136").Value = 6021.875
$ws.Range("K136").Value = 18065.625
$ws.Range("M136").Value = -15515.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3302.3333
$ws.Range("I81").Value = 3267.2778
$ws.Range("K81").Value = 6534.5556
$ws.Range("M81").Value = -5473.5556
$ws.Range("H84").Value = 3302.3333
$ws.Range("I84").Value = 3267.2778
$ws.Range("K84").Value = 32672.778
$ws.Range("M84").Value = -27368.778
$ws.Range("H96").Value = 3158.7
$ws.Range("I96").Value = 3374.75
$ws.Range("K96").Value = 3374.75
$ws.Range("M96").Value = -2001.75
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 3618.1667
$ws.Range("I132").Value = 2133.5454
$ws.Range("J132").Value = 19949
$ws.Range("K132").Value = 6400.6362
$ws.Range("L132").Value = 59847
$ws.Range("M132").Value = -3870.6362
$ws.Range("N132").Value = -64907
$ws.Range("H136").Value = 9601.429
$ws.Range("I136").Value = 10032.308
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 30096.924
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -27546.924
$ws.Range("N136").Value = -17100
